$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D (Price) and E (Volume) columns to be treated as text so that
# numeric-looking values (e.g. "1.007") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '20.478.05'
$ws.Range("E2").Value = '  -7.18%  '

$ws.Range("D3").Value = '1.449.45'
$ws.Range("E3").Value = '  -6.92%  '

$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.64%  '

$ws.Range("D5").Value = '1.007'
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("D6").Value = '278.01'
$ws.Range("E6").Value = '  -4.59%  '

$ws.Range("D7").Value = '0.3720'
$ws.Range("E7").Value = '  -5.63%  '

$ws.Range("D8").Value = '0.3094'
$ws.Range("E8").Value = '  -4.26%  '

$ws.Range("D9").Value = '40.93'
$ws.Range("E9").Value = '  -7.61%  '

$ws.Range("D10").Value = '1.013'
$ws.Range("E10").Value = '  -6.01%  '

$ws.Range("D11").Value = '0.06602'
$ws.Range("E11").Value = '  -8.81%  '

$ws.Range("E12").Value = '  +0.96%  '

$ws.Range("D13").Value = '5.434'
$ws.Range("E13").Value = '  -4.49%  '

$ws.Range("D14").Value = '17.29'
$ws.Range("E14").Value = '  -7.79%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '6.149'
$ws.Range("E15").Value = '  -7.44%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.452.42'
$ws.Range("E16").Value = '  -6.77%  '

$ws.Range("D17").Value = '0.00001021'
$ws.Range("E17").Value = '  -8.76%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '77.11'
$ws.Range("E18").Value = '  -7.81%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.05919'
$ws.Range("E19").Value = '  -10.25%  '

$ws.Range("E20").Value = '  +0.70%  '

$ws.Range("D21").Value = '5.745'
$ws.Range("E21").Value = '  -8.08%  '

$ws.Range("D22").Value = '14.55'
$ws.Range("E22").Value = '  -6.39%  '

$ws.Range("D23").Value = '10.95'
$ws.Range("E23").Value = '  -3.22%  '

$ws.Range("D24").Value = '2.310'
$ws.Range("E24").Value = '  -2.08%  '

$ws.Range("D25").Value = '20.466.46'
$ws.Range("E25").Value = '  -7.26%  '

$ws.Range("D26").Value = '2.241'
$ws.Range("E26").Value = '  -7.00%  '

$ws.Range("D27").Value = '143.17'
$ws.Range("E27").Value = '  -3.38%  '

$ws.Range("D28").Value = '17.10'
$ws.Range("E28").Value = '  -8.06%  '

$ws.Range("D29").Value = '1.620.61'
$ws.Range("E29").Value = '  -6.55%  '

$ws.Range("D30").Value = '109.45'
$ws.Range("E30").Value = '  -8.00%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '0.9218'
$ws.Range("E31").Value = '  -6.24%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '5.414'
$ws.Range("E32").Value = '  -8.50%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '3.464'
$ws.Range("E33").Value = '  -28.94%  '

$ws.Range("D34").Value = '0.07756'
$ws.Range("E34").Value = '  -6.67%  '

$ws.Range("D35").Value = '8.347'
$ws.Range("E35").Value = '  -8.91%  '

$ws.Range("D36").Value = '10.97'
$ws.Range("E36").Value = '  +1.96%  '

$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '1.007'
$ws.Range("E37").Value = '  +0.64%  '

$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = '1.415'
$ws.Range("E38").Value = '  -12.25%  '

$ws.Range("D39").Value = '0.05643'
$ws.Range("E39").Value = '  -6.30%  '

$ws.Range("D40").Value = '4.758'
$ws.Range("E40").Value = '  -7.21%  '

$ws.Range("D41").Value = '1.129'
$ws.Range("E41").Value = '  -6.39%  '

$ws.Range("D42").Value = '0.02044'
$ws.Range("E42").Value = '  -9.83%  '

$ws.Range("D43").Value = '0.1913'
$ws.Range("E43").Value = '  -7.07%  '

$ws.Range("D44").Value = '3.596'

$ws.Range("D45").Value = '0.5333'
$ws.Range("E45").Value = '  -8.34%  '

$ws.Range("E46").Value = '  -6.86%  '

$ws.Range("D47").Value = '0.5163'
$ws.Range("E47").Value = '  -7.43%  '

$ws.Range("D48").Value = '110.95'
$ws.Range("E48").Value = '  -5.83%  '

$ws.Range("E49").Value = '  -5.93%  '

$ws.Range("D50").Value = '1.062'
$ws.Range("E50").Value = '  -6.67%  '

$ws.Range("E51").Value = '  +0.68%  '

# Restore default (Normal) style on the price/volume columns so no stray
# cell formatting is introduced by the text-forcing step above.
$ws.Range("D2:E51").Style = "Normal"
